$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 525
$ws.Range("I12").Value = 450
$ws.Range("J12").Value = 600
$ws.Range("K12").Value = 450
$ws.Range("L12").Value = 600
$ws.Range("M12").Value = -280
$ws.Range("N12").Value = -940
$ws.Range("H64").Value = 3068.9656
$ws.Range("I64").Value = 2962.5
$ws.Range("J64").Value = 3200
$ws.Range("K64").Value = 2962.5
$ws.Range("L64").Value = 3200
$ws.Range("M64").Value = -2714.5
$ws.Range("N64").Value = -3696
$ws.Range("H67").Value = 3068.9656
$ws.Range("I67").Value = 2962.5
$ws.Range("J67").Value = 3200
$ws.Range("K67").Value = 2962.5
$ws.Range("L67").Value = 3200
$ws.Range("M67").Value = -2104.5
$ws.Range("N67").Value = -4916
$ws.Range("H74").Value = 2743.2856
$ws.Range("I74").Value = 2667.1667
$ws.Range("J74").Value = 3200
$ws.Range("K74").Value = 2667.1667
$ws.Range("L74").Value = 3200
$ws.Range("M74").Value = -1731.1667
$ws.Range("N74").Value = -5072
$ws.Range("H76").Value = 3007.9524
$ws.Range("I76").Value = 3018.35
$ws.Range("J76").Value = 2800
$ws.Range("K76").Value = 3018.35
$ws.Range("L76").Value = 2800
$ws.Range("M76").Value = -2703.35
$ws.Range("N76").Value = -3430
$ws.Range("H77").Value = 2743.2856
$ws.Range("I77").Value = 2667.1667
$ws.Range("J77").Value = 3200
$ws.Range("K77").Value = 13335.8335
$ws.Range("L77").Value = 16000
$ws.Range("M77").Value = -8655.833500000001
$ws.Range("N77").Value = -25360
$ws.Range("H79").Value = 3007.9524
$ws.Range("I79").Value = 3018.35
$ws.Range("J79").Value = 2800
$ws.Range("K79").Value = 3018.35
$ws.Range("L79").Value = 2800
$ws.Range("M79").Value = -1926.35
$ws.Range("N79").Value = -4984
$ws.Range("H132").Value = 17191.533
$ws.Range("I132").Value = 2325.3125
$ws.Range("J132").Value = 76656.414
$ws.Range("K132").Value = 6975.9375
$ws.Range("L132").Value = 229969.242
$ws.Range("M132").Value = -4445.9375
$ws.Range("N132").Value = -235029.242

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2922.8125
$ws.Range("I61").Value = 2288.875
$ws.Range("K61").Value = 2288.875
$ws.Range("M61").Value = -2076.875
$ws.Range("H63").Value = 2713.2
$ws.Range("I63").Value = 2244.3333
$ws.Range("J63").Value = 3416.5
$ws.Range("K63").Value = 2244.3333
$ws.Range("L63").Value = 3416.5
$ws.Range("M63").Value = -1558.3333
$ws.Range("N63").Value = -4788.5
$ws.Range("H66").Value = 2713.2
$ws.Range("I66").Value = 2244.3333
$ws.Range("J66").Value = 3416.5
$ws.Range("K66").Value = 11221.6665
$ws.Range("L66").Value = 17082.5
$ws.Range("M66").Value = -7789.666499999999
$ws.Range("N66").Value = -23946.5
$ws.Range("H74").Value = 1580.8667
$ws.Range("I74").Value = 1303.3158
$ws.Range("J74").Value = 3087.5715
$ws.Range("K74").Value = 1303.3158
$ws.Range("L74").Value = 3087.5715
$ws.Range("M74").Value = -429.3158000000001
$ws.Range("N74").Value = -4835.5715
$ws.Range("H77").Value = 1580.8667
$ws.Range("I77").Value = 1303.3158
$ws.Range("J77").Value = 3087.5715
$ws.Range("K77").Value = 6516.579000000001
$ws.Range("L77").Value = 15437.8575
$ws.Range("M77").Value = -2148.579000000001
$ws.Range("N77").Value = -24173.8575
$ws.Range("H132").Value = 13160831
$ws.Range("I132").Value = 27779880
$ws.Range("J132").Value = 3687.15
$ws.Range("K132").Value = 83339640
$ws.Range("L132").Value = 11061.45
$ws.Range("M132").Value = -83337110
$ws.Range("N132").Value = -16121.45
$ws.Range("H136").Value = 2922.8125
$ws.Range("I136").Value = 2288.875
$ws.Range("K136").Value = 6866.625
$ws.Range("M136").Value = -4316.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1442.7727
$ws.Range("I94").Value = 1411.05
$ws.Range("J94").Value = 1760
$ws.Range("K94").Value = 1411.05
$ws.Range("L94").Value = 1760
$ws.Range("M94").Value = -960.05
$ws.Range("N94").Value = -2662
$ws.Range("H105").Value = 4217.2856
$ws.Range("I105").Value = 2670
$ws.Range("J105").Value = 4639.273
$ws.Range("K105").Value = 2670
$ws.Range("L105").Value = 4639.273
$ws.Range("M105").Value = -923
$ws.Range("N105").Value = -8133.273
$ws.Range("H134").Value = 2660.724
$ws.Range("I134").Value = 2167.0454
$ws.Range("K134").Value = 6501.1362
$ws.Range("M134").Value = -3966.1362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 14287361
$ws.Range("I58").Value = 930.1
$ws.Range("J58").Value = 33335936
$ws.Range("K58").Value = 930.1
$ws.Range("L58").Value = 33335936
$ws.Range("M58").Value = -727.1
$ws.Range("N58").Value = -33336342
$ws.Range("H62").Value = 3205.8235
$ws.Range("I62").Value = 3200
$ws.Range("J62").Value = 3216.5
$ws.Range("K62").Value = 3200
$ws.Range("L62").Value = 3216.5
$ws.Range("M62").Value = -2576
$ws.Range("N62").Value = -4464.5
$ws.Range("H65").Value = 3205.8235
$ws.Range("I65").Value = 3200
$ws.Range("J65").Value = 3216.5
$ws.Range("K65").Value = 16000
$ws.Range("L65").Value = 16082.5
$ws.Range("M65").Value = -12880
$ws.Range("N65").Value = -22322.5
$ws.Range("H132").Value = 25426.05
$ws.Range("I132").Value = 1332.0426
$ws.Range("J132").Value = 112535.16
$ws.Range("K132").Value = 3996.1278
$ws.Range("L132").Value = 337605.48
$ws.Range("M132").Value = -1466.1278
$ws.Range("N132").Value = -342665.48
$ws.Range("H133").Value = 41519.582
$ws.Range("J133").Value = 41519.582
$ws.Range("L133").Value = 41519.582
$ws.Range("N133").Value = -46579.582
$ws.Range("H134").Value = 299117.9
$ws.Range("I134").Value = 1033
$ws.Range("J134").Value = 1274668.5
$ws.Range("K134").Value = 3099
$ws.Range("L134").Value = 3824005.5
$ws.Range("M134").Value = -564
$ws.Range("N134").Value = -3829075.5
$ws.Range("H136").Value = 14287361
$ws.Range("I136").Value = 930.1
$ws.Range("J136").Value = 33335936
$ws.Range("K136").Value = 2790.3
$ws.Range("L136").Value = 100007808
$ws.Range("M136").Value = -240.3000000000002
$ws.Range("N136").Value = -100012908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4317.769
$ws.Range("I34").Value = 443.66666
$ws.Range("J34").Value = 5480
$ws.Range("K34").Value = 1330.99998
$ws.Range("L34").Value = 16440
$ws.Range("M34").Value = -1246.99998
$ws.Range("N34").Value = -16608
$ws.Range("H39").Value = 2925
$ws.Range("J39").Value = 2925
$ws.Range("L39").Value = 8775
$ws.Range("N39").Value = -9363
$ws.Range("H55").Value = 2868.4211
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 2918.919
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 8756.757
$ws.Range("M55").Value = -2823
$ws.Range("N55").Value = -9110.757
$ws.Range("H68").Value = 5016350
$ws.Range("I68").Value = 6945308.5
$ws.Range("J68").Value = 4465219
$ws.Range("K68").Value = 20835925.5
$ws.Range("L68").Value = 13395657
$ws.Range("M68").Value = -20835114.5
$ws.Range("N68").Value = -13397279
$ws.Range("H71").Value = 5016350
$ws.Range("I71").Value = 6945308.5
$ws.Range("J71").Value = 4465219
$ws.Range("K71").Value = 62507776.5
$ws.Range("L71").Value = 40186971
$ws.Range("M71").Value = -62503720.5
$ws.Range("N71").Value = -40195083
$ws.Range("H131").Value = 3030.6545
$ws.Range("I131").Value = 8787.083000000001
$ws.Range("J131").Value = 1424.2094
$ws.Range("K131").Value = 26361.249
$ws.Range("L131").Value = 4272.6282
$ws.Range("M131").Value = -21321.249
$ws.Range("N131").Value = -14352.6282

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5670.732
$ws.Range("I70").Value = 5629.4116
$ws.Range("J70").Value = 5871.4287
$ws.Range("K70").Value = 5629.4116
$ws.Range("L70").Value = 5871.4287
$ws.Range("M70").Value = -5359.4116
$ws.Range("N70").Value = -6411.4287
$ws.Range("H73").Value = 5670.732
$ws.Range("I73").Value = 5629.4116
$ws.Range("J73").Value = 5871.4287
$ws.Range("K73").Value = 5629.4116
$ws.Range("L73").Value = 5871.4287
$ws.Range("M73").Value = -4693.4116
$ws.Range("N73").Value = -7743.4287
$ws.Range("H80").Value = 5484.4614
$ws.Range("I80").Value = 4687.5
$ws.Range("J80").Value = 6038.8696
$ws.Range("K80").Value = 4687.5
$ws.Range("L80").Value = 6038.8696
$ws.Range("M80").Value = -3689.5
$ws.Range("N80").Value = -8034.8696
$ws.Range("H83").Value = 5484.4614
$ws.Range("I83").Value = 4687.5
$ws.Range("J83").Value = 6038.8696
$ws.Range("K83").Value = 23437.5
$ws.Range("L83").Value = 30194.348
$ws.Range("M83").Value = -18445.5
$ws.Range("N83").Value = -40178.348
$ws.Range("H132").Value = 2940.88
$ws.Range("I132").Value = 2039.1538
$ws.Range("J132").Value = 3917.75
$ws.Range("K132").Value = 6117.4614
$ws.Range("L132").Value = 11753.25
$ws.Range("M132").Value = -3587.4614
$ws.Range("N132").Value = -16813.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3060.375
$ws.Range("I132").Value = 2197.2334
$ws.Range("J132").Value = 4498.9443
$ws.Range("K132").Value = 6591.7002
$ws.Range("L132").Value = 13496.8329
$ws.Range("M132").Value = -4061.7002
$ws.Range("N132").Value = -18556.8329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 14707876
$ws.Range("I132").Value = 1595.3182
$ws.Range("J132").Value = 41669390
$ws.Range("K132").Value = 4785.9546
$ws.Range("L132").Value = 125008170
$ws.Range("M132").Value = -2255.9546
$ws.Range("N132").Value = -125013230
$ws.Range("H136").Value = 271248.1
$ws.Range("I136").Value = 303799.47
$ws.Range("J136").Value = 2699.5
$ws.Range("K136").Value = 911398.4099999999
$ws.Range("L136").Value = 8098.5
$ws.Range("M136").Value = -908848.4099999999
$ws.Range("N136").Value = -13198.5
